$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the result value for the 4th test row (C5): record the "invalid
# login" expectation for the registration data-driven test.
$ws.Range("C5").Value = "Login(Login failed for invalid)"
